# Add cantrals by cantons
# Restructure Sheet1: merge the two header rows into a single header row,
# rename/re-order columns, and remove the now-redundant units row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------
# 1) Rewrite the first header row (currently spread across rows 1-2)
#    with the new column headers, in the new order.
# ---------------------------------------------------------------
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# A1:E1 use the plain default style (font Arial 10, General format).
$ws.Range("A1:E1").Font.Size = 10

# F1:K1 use the "Arial 9 / General" style (same font as the data rows,
# but without an explicit number-format override). Build it as a
# temporary named style so it gets its own cell-format record, then drop
# the name again so the cells keep plain (unnamed) direct formatting.
$tempStyle = $wb.Styles.Add("TempHeaderUnits", "Normal")
$tempStyle.Font.Name = "Arial"
$tempStyle.Font.Size = 9
$ws.Range("F1:K1").Style = "TempHeaderUnits"
$tempStyle.Delete()

# ---------------------------------------------------------------
# 2) Remove the old units row (row 2), which shifts the data rows
#    (previously rows 3-6) up to become rows 2-5.
# ---------------------------------------------------------------
$ws.Rows.Item(2).Delete()

# ---------------------------------------------------------------
# 3) Restore the selection state recorded after the edit.
# ---------------------------------------------------------------
$ws.Range("A2:K2").Select()
